# Add a new Job Posting row (JD_004 / Senior RPA Developer) to Sheet1.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$jobDescription = "We are seeking a Junior RPA Developer to design, develop, and support automation solutions.`nCollaborate with teams to streamline business processes using RPA tools like UiPath or Automation Anywhere. Join Akkodis to grow your skills in a dynamic, tech-driven environment"

$newRow = 5
$ws.Cells.Item($newRow, 1).Value = "JD_004"
$ws.Cells.Item($newRow, 2).Value = "Senior RPA Developer"
$ws.Cells.Item($newRow, 3).Value = $jobDescription
$ws.Cells.Item($newRow, 4).Value = 2
$ws.Cells.Item($newRow, 5).Value = 4
$ws.Cells.Item($newRow, 6).Value = 0
$ws.Cells.Item($newRow, 7).Value = 0
$ws.Cells.Item($newRow, 8).Value = 0
$ws.Cells.Item($newRow, 9).Value = 0

# The multi-line description triggers Excel's automatic row-height bump;
# AutoFit restores the row to its natural (no explicit height override) state.
$ws.Rows.Item($newRow).AutoFit()
